# edit.ps1 - applies the diff: rewrites the "Cultural Diversity" document
# into "The Orchestra of Life: An Introduction to Biology", including the
# title, author, email, body paragraphs and summary, and appends a new
# trailing empty paragraph.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: find $findText anywhere in the document and replace it with
# $replaceText. Word's Find/Replace engine merges the newly-typed text
# into whichever run immediately precedes the match when the formatting
# is identical; to keep the replacement text as its own distinct run
# (matching the source document's run layout) we toggle Bold on/off on
# the resulting range right afterwards, which forces the text-engine to
# split it back out into a dedicated run without altering its visible
# formatting.
# ---------------------------------------------------------------------
function Replace-Text {
    param(
        [string]$findText,
        [string]$replaceText
    )

    $locate = $d.Content
    $ok = $locate.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        Write-Host "NOT FOUND:" $findText
        return
    }
    $startPos = $locate.Start

    $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null

    $endPos = $startPos + $replaceText.Length
    $r = $d.Range($startPos, $endPos)
    $r.Font.Bold = 1
    $r.Font.Bold = 0
}

# ---------------- Title ----------------
Replace-Text "Cultural Diversity: Unveiling a Tapestry of Inclusion" "The Orchestra of Life: An Introduction to Biology"

# ---------------- Author (Dr. Clara Emerson -> Martha Kendrick) ----------------
Replace-Text "Dr. Clara Emerson" "Martha Kendrick"

# ---------------- Email (cemerson@researchhub.edu -> martha.kendrick@berkeley.edu) ----------------
# Change the first run's text, keep the existing "." and "edu" runs
# untouched, and insert two brand-new runs ("." and "kendrick@berkeley")
# between them.
$emailFind = $d.Content
$emailFind.Find.Execute("cemerson@researchhub", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$emailStart = $emailFind.Start
$d.Content.Find.Execute("cemerson@researchhub", $true, $false, $false, $false, $false, $true, 1, $false, "martha.kendrick@berkeley", 2) | Out-Null

$marthaLen = "martha".Length
$dotLen = ".".Length
$kendrickLen = "kendrick@berkeley".Length

$b0 = $emailStart
$b1 = $b0 + $marthaLen
$b2 = $b1 + $dotLen
$b3 = $b2 + $kendrickLen

$rMartha = $d.Range($b0, $b1)
$rMartha.Font.Bold = 1
$rMartha.Font.Bold = 0

$rDot = $d.Range($b1, $b2)
$rDot.Font.Bold = 1
$rDot.Font.Bold = 0

$rKendrick = $d.Range($b2, $b3)
$rKendrick.Font.Bold = 1
$rKendrick.Font.Bold = 0

# ---------------- Body paragraph ----------------
Replace-Text "In a world brimming with vibrancy and complexity, cultural diversity stands as a cornerstone of human existence" "Biology, the study of life, unveils the secrets of the natural world like a conductor revealing the symphony of an orchestra"

Replace-Text " Like a kaleidoscope of colors, each culture contributes its unique patterns and hues to the tapestry of our shared humanity" " We embark on a grand exploration, uncovering the intricate mechanisms of organisms from microscopic cells to towering trees"

Replace-Text " From the vibrant attire of traditional garments to the melodies of diverse musical traditions, cultural diversity weaves a rich narrative of human experiences. Delving into the depths of cultural diversity offers profound insights into the intricate relationships between identity, heritage, and the intricate tapestry of our global community." " This journey of discovery invites us to comprehend the myriad forms of life and their remarkable adaptations"

Replace-Text "As we navigate the tapestry of cultural diversity, we encounter an eclectic mix of languages, conveying a myriad of perspectives and narratives" "Imagine a vast tapestry of ecosystems, each interwoven with diverse organisms, interacting like instruments in a grand symphony"

Replace-Text " Each language embodies a distinct worldview, shaping the ways in which its speakers perceive and interact with the world around them" " Ecosystems, like perfectly orchestrated compositions, are intricate webs of relationships where every organism plays a unique role"

Replace-Text " The nuances of linguistic expressions reveal cultural values, beliefs, and norms, providing a glimpse into the hearts and minds of diverse communities. Moreover, language serves as a bridge, connecting people from different cultural backgrounds, fostering understanding, and promoting harmonious coexistence." " With the zeal of an aspiring musician, we shall delve into these ecosystems to decipher this harmonious interplay of life"

Replace-Text "Cultural diversity extends beyond language and encompasses a multitude of expressions, including art, music, dance, and cuisine" "Finally, we unravel the molecular dance of life, understanding the intricacies of DNA, the blueprint of living beings"

Replace-Text " These creative endeavors reflect the soul of a culture, showcasing its history, traditions, and aspirations" " This microscopic blueprint holds the secrets to genetic inheritance and evolution, like a beautifully orchestrated code that ensures the continuity of life itself"

Replace-Text " In the vibrant strokes of traditional paintings, we witness the stories of ancestors, while the melodies of folk songs transport us to distant lands, whispering tales of love, loss, and triumph. Through the rhythmic movements of traditional dances, we connect with the spirit of a culture, feeling the pulse of its heartbeat. And as we savor the delectable flavors of diverse cuisines, we embark on a culinary journey, tasting the essence of distant lands and forging bonds of shared experience." " Biology's captivating journey through the molecular realm unveils the very essence of existence"

# ---------------- Summary paragraph ----------------
Replace-Text "Cultural diversity is a testament to the extraordinary richness and complexity of human existence" "Biology, an enthralling adventure into life's depths, unlocks the wonders of ecosystems, organisms, and the molecular foundation of life"

Replace-Text " It encompasses a multitude of expressions, including language, art, music, dance, and cuisine, each contributing to the tapestry of our shared humanity. Embracing cultural diversity fosters understanding, promotes harmonious coexistence, and enriches our lives with a kaleidoscope of colors, melodies, and flavors. As we delve deeper into the nuances of different cultures, we gain a profound appreciation for the interconnectedness of our global community, recognizing that our differences are sources of strength and unity rather than division." " From the grandeur of ecosystems to the intricate structures of cells, from the marvel of genetic inheritance to the captivating molecular dance within, biology unveils a symphony that inspires awe and appreciation for the beauty and complexity of life"

# ---------------- Append trailing empty paragraph ----------------
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
